# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.333.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.82%  "

$ws.Range("D3").Value = "'2.426.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.00%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'573.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").Value = "'142.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.20%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "'0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.67%  "

$ws.Range("D9").Value = "'2.427.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.07%  "

$ws.Range("D10").Value = "'0.106"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.29%  "

$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("E12").Value = "  -1.66%  "

$ws.Range("D13").Value = "'0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.53%  "

$ws.Range("D14").Value = "'26.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.49%  "

$ws.Range("D15").Value = "'0.0000172"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.37%  "

$ws.Range("D16").Value = "'2.893.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("D17").Value = "'62.312.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").Value = "'2.428.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.63%  "

$ws.Range("D19").Value = "'10.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.77%  "

$ws.Range("D20").Value = "'7.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.78%  "

$ws.Range("D21").Value = "'328.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("D22").Value = "'4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "

$ws.Range("D23").Value = "'1.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.09%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "'65.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("D26").Value = "'626.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.59%  "

$ws.Range("D27").Value = "'8.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.98%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'2.563.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0₃0956"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.93%  "

$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("D31").Value = "'1.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.10%  "

$ws.Range("D32").Value = "'8.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.20%  "

$ws.Range("D33").Value = "'1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.89%  "

$ws.Range("E34").Value = "  -4.26%  "

$ws.Range("D35").Value = "'4.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.50%  "

$ws.Range("E36").Value = "  +0.47%  "

$ws.Range("D37").Value = "'1.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.75%  "

$ws.Range("D38").Value = "'0.373"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.38%  "

$ws.Range("D39").Value = "'18.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.44%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'147.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").Value = "'5.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.80%  "

$ws.Range("D42").Value = "'1.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.45%  "

$ws.Range("D43").Value = "'42.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "'2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.61%  "

$ws.Range("D46").Value = "'143.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.46%  "

$ws.Range("D47").Value = "'3.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.26%  "

$ws.Range("D48").Value = "'0.0520"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.22%  "

$ws.Range("D49").Value = "'0.596"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("D50").Value = "'19.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.37%  "

$ws.Range("D51").Value = "'0.0₆0231"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.10%  "
